# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (same layout as the other quarterly
#    fund-holding sheets) right before the "总计" summary sheet.
# 2. Add a new top row to the "总计" summary sheet for the 2022-Q1 quarter
#    and renumber the running index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, positioned before "总计"
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$q1 = $wb.Worksheets.Add($totalSheet, $null)
$q1.Name = "2022-Q1"

# Bring over the header row + index-column styling from an existing
# quarterly sheet so fonts/borders/alignment match exactly.
$template.Range("B1:H1").Copy($q1.Range("B1:H1"))
$template.Range("A2:A3").Copy($q1.Range("A2:A3"))

# Fund rows for 2022-Q1. Columns B and D:G hold numeric-looking text in
# the source data ("012588", "37.54", ...), so force a text format
# before writing, then drop back to the default style (matches how the
# other quarterly sheets store these values: plain cells, no explicit
# style, shared-string text).
$q1.Range("B2:B3").NumberFormat = "@"
$q1.Range("D2:G3").NumberFormat = "@"

$q1.Range("B2").Value = "012588"
$q1.Range("C2").Value = "南方港股通优势企业混合型证券投资基金A"
$q1.Range("D2").Value = "37.54"
$q1.Range("E2").Value = "71.00"
$q1.Range("F2").Value = "2.78"
$q1.Range("G2").Value = "1.0436"
$q1.Range("H2").Value = 5

$q1.Range("B3").Value = "012589"
$q1.Range("C3").Value = "南方港股通优势企业混合型证券投资基金C"
$q1.Range("D3").Value = "2.05"
$q1.Range("E3").Value = "71.00"
$q1.Range("F3").Value = "2.78"
$q1.Range("G3").Value = "0.0570"
$q1.Range("H3").Value = 5

$q1.Range("B2:B3").Style = "Normal"
$q1.Range("D2:G3").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet with the new quarter
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.1

# Re-copy the index-column style onto the freshly inserted row, then
# renumber the running 0-based index for every data row.
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
